# Regenerate the "K" (strikeouts) column values for Ryan Pressly's 2022
# save_data sheet. The commit message indicates the G column was switched
# from a "Strike#" style stat to "K", so the per-row values are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:G62, in row order (row 2 first).
$newValues = @(0,2,1,1,1,0,3,3,1,1,2,2,2,1,1,1,1,2,2,3,2,1,0,2,0,3,2,2,2,2,0,3,3,2,1,3,2,1,2,0,2,1,0,0,2,0,1,0,1,1,0,1,3,0,1,0,0,1,1,0,1)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 7).Value = $val   # column G = 7
    $row++
}
